$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.188.30"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "3.047.37"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("E4").Value = "  -0.10%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "590.74"
$ws.Range("E5").Value = "  +0.35%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "152.39"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("E7").Value = "  +0.02%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.541"
$ws.Range("E8").Value = "  -1.40%  "

$ws.Range("D9").Value = "3.045.91"
$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("E12").Value = "  -2.50%  "

$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "36.52"
$ws.Range("E13").Value = "  -1.81%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.0000236"
$ws.Range("E14").Value = "  -2.23%  "

$ws.Range("E15").Value = "  +1.56%  "

$ws.Range("D16").Value = "3.551.98"
$ws.Range("E16").Value = "  -0.78%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "7.16"
$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("D18").Value = "63.156.85"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").Value = "3.049.43"
$ws.Range("E19").Value = "  -0.74%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "481.05"
$ws.Range("E20").Value = "  +1.07%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.27"
$ws.Range("E21").Value = "  -2.42%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.706"
$ws.Range("E22").Value = "  -1.70%  "

$ws.Range("E23").Value = "  -0.08%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.43"
$ws.Range("E24").Value = "  +2.22%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "82.12"
$ws.Range("E25").Value = "  +0.81%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "12.75"
$ws.Range("E26").Value = "  -2.08%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.70"
$ws.Range("E27").Value = "  +8.39%  "

$ws.Range("E28").Value = "  +0.13%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.39"
$ws.Range("E29").Value = "  +1.30%  "

$ws.Range("E30").Value = "  +0.15%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.21"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("E32").Value = "  -0.09%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "27.51"
$ws.Range("E33").Value = "  +0.84%  "

$ws.Range("E34").Value = "  -2.98%  "

$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0820"
$ws.Range("E35").Value = "  -3.22%  "

$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.06"
$ws.Range("E36").Value = "  +0.38%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "3.29"
$ws.Range("E37").Value = "  -1.37%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "5.93"
$ws.Range("E38").Value = "  -3.07%  "

$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("E40").Value = "  -0.57%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "50.39"
$ws.Range("E41").Value = "  +0.19%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "436.46"
$ws.Range("E42").Value = "  -1.55%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.289"
$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("E44").Value = "  +3.08%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0363"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "2.828.38"
$ws.Range("E46").Value = "  +0.75%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "38.25"
$ws.Range("E47").Value = "  -4.46%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "130.16"
$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("E49").Value = "  +0.01%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "25.16"
$ws.Range("E50").Value = "  +0.15%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.22"
$ws.Range("E51").Value = "  -1.80%  "
